# Auto-generated edit script applying the diff to 杭州-漫展信息.xlsx
# Updates "想去人数" (F column) counters and one venue address/image change (row 20)
# across sheets 展览, 演出, 本地生活, 全部类型, matching commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 397
$ws.Range("F3").Value = 1082
$ws.Range("F4").Value = 9642
$ws.Range("F5").Value = 207
$ws.Range("F6").Value = 76
$ws.Range("F8").Value = 6574
$ws.Range("F10").Value = 10473
$ws.Range("F11").Value = 11617
$ws.Range("F12").Value = 1257
$ws.Range("F13").Value = 1202
$ws.Range("F14").Value = 5036
$ws.Range("F15").Value = 833
$ws.Range("F16").Value = 496
$ws.Range("F18").Value = 343
$ws.Range("F19").Value = 183
$ws.Range("D20").Value = "创意路1号 中国智谷富春园区"
$ws.Range("I20").Value = "//i2.hdslb.com/bfs/openplatform/202407/e9atNS5Y1720172990106.png"
$ws.Range("F21").Value = 278
$ws.Range("F22").Value = 1907
$ws.Range("F24").Value = 1319
$ws.Range("F27").Value = 2079
$ws.Range("F29").Value = 665
$ws.Range("F30").Value = 2743
$ws.Range("F32").Value = 1835
$ws.Range("F34").Value = 831
$ws.Range("F35").Value = 86
$ws.Range("F36").Value = 936
$ws.Range("F37").Value = 25
$ws.Range("F38").Value = 55
$ws.Range("F39").Value = 3431
$ws.Range("F42").Value = 532
$ws.Range("F43").Value = 596
$ws.Range("F46").Value = 253
$ws.Range("F47").Value = 10
$ws.Range("F48").Value = 4237
$ws.Range("F49").Value = 83

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 32
$ws.Range("F26").Value = 45
$ws.Range("F29").Value = 2

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6094

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 397
$ws.Range("F3").Value = 1082
$ws.Range("F4").Value = 9643
$ws.Range("F5").Value = 76
$ws.Range("F9").Value = 10473
$ws.Range("F10").Value = 11617
$ws.Range("F12").Value = 1202
$ws.Range("F13").Value = 5036
$ws.Range("F14").Value = 833
$ws.Range("F15").Value = 496
$ws.Range("F17").Value = 343
$ws.Range("F18").Value = 32
$ws.Range("F19").Value = 183
$ws.Range("D20").Value = "创意路1号 中国智谷富春园区"
$ws.Range("I20").Value = "//i2.hdslb.com/bfs/openplatform/202407/e9atNS5Y1720172990106.png"
$ws.Range("F21").Value = 278
$ws.Range("F22").Value = 1907
$ws.Range("F24").Value = 1319
$ws.Range("F26").Value = 2079
$ws.Range("F28").Value = 665
$ws.Range("F29").Value = 2743
$ws.Range("F31").Value = 1835
$ws.Range("F34").Value = 831
$ws.Range("F38").Value = 86
$ws.Range("F39").Value = 936
$ws.Range("F40").Value = 25
$ws.Range("F41").Value = 45
$ws.Range("F44").Value = 532
$ws.Range("F45").Value = 596
$ws.Range("F47").Value = 253
$ws.Range("F48").Value = 10
$ws.Range("F49").Value = 4237

